# ---------------------------------------------------------------------------
# Edit 1: merge the three runs that make up
#   "    " + "School Band Vice Captain 2019-20" + ":"
# into a single run "    School Band Vice Captain 2019-20:" while keeping
# the formatting (rPr) and rsid of the very first run.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*School Band Vice Captain 2019-20:*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    # Keep the leading run (the 4 preserved spaces) untouched; remove the
    # two runs that follow it ("School Band Vice Captain 2019-20" and ":")
    # and re-insert their combined text straight after the first run so it
    # gets folded into that run instead of creating new ones.
    $tail = $d.Range($pStart + 4, $pEnd - 1)
    $tailText = $tail.Text
    $tail.Delete()

    $head = $d.Range($pStart, $pStart + 4)
    $head.InsertAfter($tailText)
}

# ---------------------------------------------------------------------------
# Edit 2: add the built-in "FollowedHyperlink" character style to the
# style sheet (mirrors what Word adds the first time a followed hyperlink
# is encountered).
# ---------------------------------------------------------------------------
# wdStyleTypeCharacter = 2
$followed = $d.Styles.Add("FollowedHyperlink", 2)
$followed.BaseStyle = "DefaultParagraphFont"
$followed.Priority = 99
$followed.UnhideWhenUsed = $true
$followed.Font.Underline = 1
$followed.Font.TextColor.ObjectThemeColor = 11
